$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the values that are missing from the original sheet: the new
#    TotalConfirmedNewCases (G) / TotalNewDeaths (I) columns on the existing
#    rows, the corrected header text, and the whole new row 6.
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 0

$ws.Range("A6").Value = 71
$ws.Range("B6").Value = 235
$ws.Range("C6").Value = "LATIN AMER. & CARIB    "
$ws.Range("D6").Value = 43921
$ws.Range("E6").Value = "British Virgin Islands"
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "Imported cases only"
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 5385

# ---------------------------------------------------------------------------
# 2. Build the two new cell formats (center/center alignment, one keeping the
#    default General number format, one using the new "yyyy-mm-dd;" custom
#    date format) on scratch cells far away from the used range, then copy
#    those formats onto the real ranges with PasteSpecial so every cell in
#    the destination lands on the SAME style index in one shot - applying
#    HorizontalAlignment/VerticalAlignment directly on a multi-cell range
#    leaves orphaned intermediate styles behind in the saved styles table.
# ---------------------------------------------------------------------------
$scratchGeneral = $ws.Range("ZZ1")
$scratchGeneral.HorizontalAlignment = -4108
$scratchGeneral.VerticalAlignment = -4108

$scratchDate = $ws.Range("ZZ2")
$scratchDate.NumberFormat = "yyyy-mm-dd;"
$scratchDate.HorizontalAlignment = -4108
$scratchDate.VerticalAlignment = -4108

# Every column except D (Date) gets the centered General style.
$scratchGeneral.Copy()
$ws.Range("A1:C6").PasteSpecial(-4122)
$ws.Range("E1:O6").PasteSpecial(-4122)

# Column D (Date) gets the centered custom date style.
$scratchDate.Copy()
$ws.Range("D1:D6").PasteSpecial(-4122)

$scratchGeneral.Clear()
$scratchDate.Clear()

# ---------------------------------------------------------------------------
# 3. Column widths: every used column (A:O) becomes 27 characters wide.
# ---------------------------------------------------------------------------
$ws.Range("A1:O1").ColumnWidth = 26.14
